$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.75"
$ws.Range("E2").Value = "'7.34%"
$ws.Range("D3").Value = "'32.08"
$ws.Range("E3").Value = "'9.83%"
$ws.Range("D4").Value = "'5.346"
$ws.Range("E4").Value = "'5.26%"
$ws.Range("D5").Value = "'0.07571"
$ws.Range("E5").Value = "'13.35%"
$ws.Range("D6").Value = "'7.815"
$ws.Range("E6").Value = "'6.58%"
$ws.Range("D7").Value = "'3.673"
$ws.Range("E7").Value = "'7.88%"
$ws.Range("D8").Value = "'1.580"
$ws.Range("E8").Value = "'17.07%"
$ws.Range("D9").Value = "'0.9113"
$ws.Range("E9").Value = "'-0.76%"
$ws.Range("D10").Value = "'0.01695"
$ws.Range("E10").Value = "'2,519.00%"
$ws.Range("D11").Value = "'0.1689"
$ws.Range("E11").Value = "'6.40%"
$ws.Range("D12").Value = "'0.07694"
$ws.Range("E12").Value = "'12.88%"
$ws.Range("D13").Value = "'0.08091"
$ws.Range("E13").Value = "'5.50%"
$ws.Range("D14").Value = "'0.03021"
$ws.Range("E14").Value = "'2.93%"
$ws.Range("D15").Value = "'0.09875"
$ws.Range("E15").Value = "'9.87%"
$ws.Range("D16").Value = "'0.001523"
$ws.Range("E16").Value = "'-4.05%"
$ws.Range("D17").Value = "'0.04549"
$ws.Range("E17").Value = "'1.15%"
$ws.Range("D18").Value = "'0.006526"
$ws.Range("E18").Value = "'4.45%"
$ws.Range("D19").Value = "'3.507"
$ws.Range("E19").Value = "'1.60%"
$ws.Range("D20").Value = "'2.238"
$ws.Range("E20").Value = "'0.89%"
$ws.Range("E21").Value = "'1.68%"
$ws.Range("E22").Value = "'1.98%"
$ws.Range("D23").Value = "'4.187"
$ws.Range("E23").Value = "'3.23%"
$ws.Range("E24").Value = "'2.88%"
$ws.Range("D25").Value = "'0.001216"
$ws.Range("E25").Value = "'2.08%"
$ws.Range("D26").Value = "'0.004490"
$ws.Range("E26").Value = "'9.08%"
$ws.Range("E27").Value = "'8.33%"
$ws.Range("E28").Value = "'7.57%"
$ws.Range("E40").Value = "'7.68%"
$ws.Range("D41").Value = "'0.007130"
$ws.Range("E41").Value = "'6.12%"
$ws.Range("D42").Value = "'0.1364"
$ws.Range("E42").Value = "'10.15%"
$ws.Range("E43").Value = "'8.02%"
$ws.Range("D44").Value = "'0.01393"
$ws.Range("E44").Value = "'4.08%"
$ws.Range("D45").Value = "'0.00006182"
$ws.Range("E45").Value = "'8.19%"
$ws.Range("E46").Value = "'-3.83%"
$ws.Range("E47").Value = "'-0.56%"
